$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.211629152297974
$ws.Range("B1").Value = 2.267055034637451
$ws.Range("C1").Value = 6.194401264190674
$ws.Range("D1").Value = 2.009319543838501
$ws.Range("E1").Value = 1.168142437934875
